# This workbook is a weekly rolling price log (newest date near the top).
# A new week's observation was inserted as the new row 32, pushing every
# row that was previously at 32..71 down by one (to 33..72), which is why
# the sheet's used range grows from A1:R71 to A1:R72.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at 32; everything below (old rows 32-71) shifts down
# to rows 33-72, carrying its original values/formatting with it.
$ws.Rows.Item(32).Insert()

# Populate the newly inserted row 32 with this week's record.
$ws.Range("A32").Value = 11
$ws.Range("B32").Value = "Vega Monumental Concepción"
$ws.Range("C32").Value = "Bíobío"
$ws.Range("D32").Value = 45225
$ws.Range("E32").Value = 8
$ws.Range("F32").Value = 100112026
$ws.Range("G32").Value = "Haba"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 90
$ws.Range("K32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 10000
$ws.Range("N32").Value = "$/saco 25 kilos"
$ws.Range("O32").Value = "Provincia de Melipilla"
$ws.Range("P32").Value = 400
$ws.Range("Q32").Value = 25
$ws.Range("R32").Value = "Hortaliza"
